# This script applies a data re-ordering edit to the "artfynd" worksheet.
# Several rows had their observation data shuffled (IDs, coordinates, times,
# comments, species info, etc. moved between rows) while keeping the
# constant/common columns (locality, municipality, dates, observer, ...)
# untouched. We implement this by swapping/rotating the values of the
# columns that actually change, row by row, exactly as described by the
# diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Get-RowData($row) {
    $data = @{}
    $data['A']  = $ws.Range("A$row").Value()
    $data['B']  = $ws.Range("B$row").Value()
    $data['E']  = $ws.Range("E$row").Value()
    $data['F']  = $ws.Range("F$row").Value()
    $data['G']  = $ws.Range("G$row").Value()
    $data['H']  = $ws.Range("H$row").Value()
    $data['M']  = $ws.Range("M$row").Value()
    $data['Q']  = $ws.Range("Q$row").Value()
    $data['R']  = $ws.Range("R$row").Value()
    $data['Z']  = $ws.Range("Z$row").Value()
    $data['AB'] = $ws.Range("AB$row").Value()
    $data['AC'] = $ws.Range("AC$row").Value()
    return $data
}

function Set-RowData($row, $data) {
    $ws.Range("A$row").Value = $data['A']
    $ws.Range("B$row").Value = $data['B']
    $ws.Range("E$row").Value = $data['E']
    $ws.Range("F$row").Value = $data['F']
    $ws.Range("G$row").Value = $data['G']

    if ($data['H'] -eq $null) {
        $ws.Range("H$row").ClearContents()
    } else {
        $ws.Range("H$row").Value = $data['H']
    }

    if ($data['M'] -eq $null) {
        $ws.Range("M$row").ClearContents()
    } else {
        $ws.Range("M$row").Value = $data['M']
    }

    $ws.Range("Q$row").Value  = $data['Q']
    $ws.Range("R$row").Value  = $data['R']
    $ws.Range("Z$row").Value  = $data['Z']
    $ws.Range("AB$row").Value = $data['AB']
    $ws.Range("AC$row").Value = $data['AC']
}

# --- Rows 11 and 13: swap contents ---
$row11 = Get-RowData 11
$row13 = Get-RowData 13
Set-RowData 11 $row13
Set-RowData 13 $row11

# --- Rows 12 and 14: swap contents ---
$row12 = Get-RowData 12
$row14 = Get-RowData 14
Set-RowData 12 $row14
Set-RowData 14 $row12

# --- Rows 27, 28, 29: cyclic rotation (27<-28, 28<-29, 29<-27) ---
$row27 = Get-RowData 27
$row28 = Get-RowData 28
$row29 = Get-RowData 29
Set-RowData 27 $row28
Set-RowData 28 $row29
Set-RowData 29 $row27

# --- Rows 35 and 36: swap contents ---
$row35 = Get-RowData 35
$row36 = Get-RowData 36
Set-RowData 35 $row36
Set-RowData 36 $row35
